$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2 = -2.5
    3 = 2.78
    4 = 2.78
    5 = -0.06
    6 = 0.34
    7 = 3.7
    8 = 1.28
    9 = 0.43
    10 = 0.45
    11 = 2.22
    12 = 1.35
    13 = 1.23
    14 = 2.77
    15 = 2.03
    16 = 0.62
    17 = 2.69
    18 = 2.69
    19 = 2.54
    20 = 2.31
    21 = 1.58
    22 = 2.61
    23 = 2.72
    24 = 2.12
    25 = 2.89
    26 = 3.51
    27 = 2.42
    28 = 3.29
    29 = 3.21
    30 = 3.06
    31 = 3.44
    32 = 3.27
    33 = 3.28
    34 = 3.22
    35 = 3.09
    36 = 3.02
    37 = 3.39
    38 = 2.52
    39 = 3.16
    40 = 3.61
    41 = 3.83
    42 = 2.39
    43 = 3.8
    44 = 2.56
    45 = 3.29
    46 = 3.96
    47 = 3.97
    48 = 4.67
    49 = 3.3
    50 = 3.45
    51 = 3.78
    52 = 3.54
    53 = 3.95
    54 = 3.56
    55 = 4.07
    56 = 2.37
    57 = 4.46
    58 = 4.58
    59 = 3.64
    60 = 3.71
    61 = 5.39
    62 = 4.85
    63 = 3.78
    64 = 3.97
    65 = 4.38
    66 = 5.33
    67 = 4.26
    68 = 4.63
    69 = 4.2
    70 = 3.46
    71 = 4.81
    72 = 4.08
    73 = 5
    74 = 4.28
    75 = 4.78
    76 = 5.24
    77 = 4.58
    78 = 4.46
    79 = 4.29
    80 = 4.77
    81 = 5.23
    82 = 5.16
    83 = 5.8
    84 = 4.42
    85 = 5.48
    86 = 5.09
    87 = 5.49
    88 = 6
    89 = 5.5
    90 = 4.07
    91 = 3.88
    92 = 5.87
    93 = 5.68
    94 = 6.12
    95 = 5.97
    97 = 6.24
    98 = 6.3
    99 = 7.37
    100 = 8.72
    101 = 8.07
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $updates[$row]
}

Write-Output "Updated $($updates.Count) cells in column C"
